$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) column C contains a date serial number for every
# data row (rows 2 through 98). This value is bumped forward by one day
# (46075 -> 46076), i.e. from 2026-02-22 to 2026-02-23.
$range = $ws.Range("C2:C98")
$range.Value = 46076
